$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated average_county_temperature (column I) using NOAA data, plus the
# dependent worst_ashp_cop (column N) and best_ashp_cop (column O) values
# for the affected facility rows.

# AdvanSix Resins & Chemicals LLC (PA) - rows 2,3
$ws.Range("I2").Value = 21.79166666666666
$ws.Range("N2").Value = 1.105721877767936
$ws.Range("O2").Value = 1.143718778908418
$ws.Range("I3").Value = 21.79166666666666

# Kraton Chemical, LLC (GA) - rows 10,11
$ws.Range("I10").Value = -1.819444444444444
$ws.Range("N10").Value = 1.004851086664878
$ws.Range("O10").Value = 1.035188389617639
$ws.Range("I11").Value = -1.819444444444444

# ALTIVIA Petrochemicals, LLC - HAVERHILL COMPLEX (OH) - rows 18,19
$ws.Range("I18").Value = 21.28240740740739
$ws.Range("N18").Value = 1.103333005990376
$ws.Range("O18").Value = 1.14113834478515
$ws.Range("I19").Value = 21.28240740740739

# Flint Hills Resources Joliet, LLC (IL) - rows 22,23
$ws.Range("I22").Value = 12.93898809523811
$ws.Range("N22").Value = 1.065614691876665
$ws.Range("O22").Value = 1.100460934966844
$ws.Range("I23").Value = 12.93898809523811
